$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell A1 from "SL.No." to "Q.No."
$ws.Range("A1").Value = "Q.No."

# Update the serial number column to the new "QnA_xx" labels
$ws.Range("A2").Value = "QnA_01"
$ws.Range("A3").Value = "QnA_02"
$ws.Range("A4").Value = "QnA_03"

# Add the fourth question as a new row
$ws.Range("A5").Value = "QnA_04"
$ws.Range("B5").Value = "Write a query to find the employee with the highest salary in each department in the year 2025."

# Move the active selection to C2, matching the saved view state
$ws.Range("C2").Select() | Out-Null
